# Apply the committed changes to the workbook:
#  - Beam sheet: column D (Ix) for rows 41-53 changed from 9.48E-5 to 10000
#  - Story_shear sheet: three new columns added (omega_1, omega_2, floor_area)
#    with header row + three data rows of sample values
#  - Selections/active views updated to match where the author was last
#    working (Story_shear I7, then back to Beam D41:D53 as the final active
#    sheet/selection, matching the saved workbook state)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Beam sheet — bulk-update column D (Ix) for rows 41 through 53
# ---------------------------------------------------------------
$wsBeam = $wb.Worksheets.Item("Beam")
$wsBeam.Range("D41:D53").Value = 10000

# ---------------------------------------------------------------
# 2) Story_shear sheet — add omega_1 / omega_2 / floor_area columns
# ---------------------------------------------------------------
$wsStory = $wb.Worksheets.Item("Story_shear")

$wsStory.Range("E1").Value = "omega_1"
$wsStory.Range("F1").Value = "omega_2"
$wsStory.Range("G1").Value = "floor_area"

$wsStory.Range("E2:E4").Value = 100
$wsStory.Range("F2:F4").Value = 200
$wsStory.Range("G2:G4").Value = 150

# ---------------------------------------------------------------
# 3) Restore view/selection state: visit Story_shear (last touched
#    there at I7), then return to Beam, selecting D41:D53 so Beam
#    ends up as the active sheet/selection, matching the workbook.
# ---------------------------------------------------------------
[void]$wsStory.Activate()
[void]$wsStory.Range("I7").Select()

[void]$wsBeam.Activate()
$excel.ActiveWindow.ScrollRow = 37
[void]$wsBeam.Range("D41:D53").Select()
